# Customer Group RFC 2022 Spring-Summer.xlsx
# Update the sample RFC volume values on row 2 of the "Template" sheet:
#   B2: "26" -> "09"   (Retail RFC volume)
#   C2: "40" -> "97"   (Wholesale RFC volume)
#
# Both cells must remain TEXT (shared-string) cells with no explicit
# style override, matching every other data cell in row 2 (A2/B2/C2 carry
# no "s" attribute in the saved XML). Writing a numeric-looking string
# straight into .Value gets auto-coerced into a number by Excel, so
# instead we enter a text formula (="09") that evaluates to the literal
# string, then copy/paste-special as values. That collapses the formula
# down to a plain cell while preserving the text result type - without
# ever touching NumberFormat/Style, so no stray style gets created.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

$xlPasteValues = -4163

$ws.Cells.Item(2, 2).Formula = "=""09"""
$ws.Cells.Item(2, 2).Copy()
$ws.Cells.Item(2, 2).PasteSpecial($xlPasteValues)

$ws.Cells.Item(2, 3).Formula = "=""97"""
$ws.Cells.Item(2, 3).Copy()
$ws.Cells.Item(2, 3).PasteSpecial($xlPasteValues)
